$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number (45179 = 2023-09-10)
# that needs to be bumped by one day (45180 = 2023-09-11) for every
# data row (rows 2 through 527).
$ws.Range("C2:C527").Value = 45180
